# This script reproduces the author's edit:
#  1. On sheet "Prix Spot", a new column is inserted before column DW (column
#     index 127). This shifts the existing "01-oct." .. "31-oct." block (and
#     everything after it) one column to the right, turning the old
#     DW:FA range into DX:FB. The newly inserted DW column is populated with
#     the header "22-nov" (row 1) and "-" placeholders for all the data rows
#     (rows 2-25), exactly like the neighbouring "no data yet" columns.
#  2. On sheets "Gaz" and "CO2", a new trailing row (156) is appended with the
#     date "2025-11-20" and the corresponding last-price value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Prix Spot" sheet: insert the new "22-nov" column
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Column 127 is "DW". Inserting there shifts DW:FA -> DX:FB and leaves a
# blank DW column ready to be filled in.
$wsPrix.Columns.Item(127).Insert()

# Header for the newly inserted column.
$wsPrix.Range("DW1").Value = "22-nov"

# Data rows 2-25 get the "no data" placeholder, just like the other
# still-empty future date columns.
$wsPrix.Range("DW2:DW25").Value = "-"

# ---------------------------------------------------------------------------
# 2) "Gaz" and "CO2" sheets: append the new daily price row
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A156").Value = "'2025-11-20"
$wsGaz.Range("A156").Style = "Normal"
$wsGaz.Range("B156").Value = 30.04

$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A156").Value = "'2025-11-20"
$wsCO2.Range("A156").Style = "Normal"
$wsCO2.Range("B156").Value = 80.92
